$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row to append after the last used row (row 83 -> 84)
$newRow = 84

# A84: date value 45884 (serial date), formatted same as column A (style copied from A83)
$ws.Cells.Item($newRow, 1).Value = 45884
$ws.Cells.Item($newRow, 2).Value = 0.325608361860148

# Copy the style of the cell above (A83) into the new date cell (A84) so it keeps
# the same date number format / font / border as the rest of the column.
$ws.Range("A83").Copy() | Out-Null
$ws.Range("A84").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set the value after paste (PasteSpecial formats only, so value should be intact,
# but ensure it explicitly to avoid any formula/format side effects)
$ws.Cells.Item($newRow, 1).Value = 45884
$ws.Cells.Item($newRow, 2).Value = 0.325608361860148
